# Update the password for Lambros.Poullais (Business / Manufacturer /
# AuthorisedRep rows) on every sheet so it matches the other users'
# "MHRA12345" password, then leave the selection/active-tab state the
# way the author left it after making the edits: DeviceSetupLogins and
# InjectSpecificUser both show the just-edited B14:B16 block selected,
# while Sheet1 (now the active tab) is left on C25.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Sheet1", "DeviceSetupLogins", "InjectSpecificUser")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B14").Value = "MHRA12345"
    $ws.Range("B15").Value = "MHRA12345"
    $ws.Range("B16").Value = "MHRA12345"
}

# Select the edited rows on DeviceSetupLogins, then InjectSpecificUser,
# finishing on Sheet1 (which becomes the active tab) at C25 - matching
# the end-state captured in the workbook's view settings.
$wb.Worksheets.Item("DeviceSetupLogins").Range("B14:B16").Select()
$wb.Worksheets.Item("InjectSpecificUser").Range("B14:B16").Select()
$wb.Worksheets.Item("Sheet1").Range("C25").Select()
